# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps to reflect the new report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date (shared with de-de!H2)
$wsOverview.Range("G2").Value = "2016-08-12 05:07:56"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-12 05:07:51"
$wsZhCn.Range("K2").Value = "2016-08-12 05:08:13"

# de-de sheet: Correspond Handoff Datetime (same value as Overview!G2) / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-08-12 05:07:56"
$wsDeDe.Range("K2").Value = "2016-08-12 05:08:20"
